# Automation Script for invalid login
# Adds a new "InvalidLogin" worksheet (after the existing "ValidLogin" sheet)
# containing a negative-test data table, and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Existing sheet ("ValidLogin") stays in place; the new sheet is inserted
# right after it so sheet order becomes [ValidLogin, InvalidLogin].
$validSheet = $wb.Worksheets.Item(1)

$invalidSheet = $wb.Worksheets.Add($null, $validSheet)
$invalidSheet.Name = "InvalidLogin"

# Header row
$invalidSheet.Range("A1").Value = "UserName"
$invalidSheet.Range("B1").Value = "Password"
$invalidSheet.Range("C1").Value = "ErrMsg"

# Negative-test data row
$invalidSheet.Range("A2").Value = "abcd"
$invalidSheet.Range("B2").Value = "xyz"
$invalidSheet.Range("C2").Value = "Invalid username and/or password."

# Size columns to their content, like the authored workbook.
$invalidSheet.Range("A:C").EntireColumn.AutoFit()

# Select the whole sheet (mirrors the recorded selection state) and make
# this new sheet the active tab.
[void]$invalidSheet.Cells.Select()
[void]$invalidSheet.Activate()
